$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column AG (33) - shifts AG..AJ to AH..AK
$ws.Range("AG1").EntireColumn.Insert()

# Match the newly inserted column's width to its left neighbour (AF)
$ws.Range("AG1").ColumnWidth = $ws.Range("AF1").ColumnWidth

# Set the header text for the newly inserted column (row 7)
$ws.Range("AG7").Value = "Giờ bắt đầu thổi"

# Update the view state to match target (scrolled/selected)
$excel.ActiveWindow.ScrollColumn = 17
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("AG10").Select()
